$d = $word.ActiveDocument

# Pull the whole package as flat OPC WordprocessingML so we can retarget the
# `name` attribute on the wp:docPr / pic:cNvPr pair for each inline picture.
# (InlineShape has no scriptable "Name" property in the Word object model --
# only AlternativeText/Title -- so the rename has to go through the raw XML.)
$x = $d.WordOpenXML

# --- Footer 1 (docPr id="3"): Pearson logo image1.png -> image2.png ---
$old1 = 'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"'
$new1 = 'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"'
$x = $x.Replace($old1, $new1)

# --- Footer 2 (docPr id="2"): Pearson logo image1.png -> image2.png ---
$old2 = 'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"'
$new2 = 'wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"'
$x = $x.Replace($old2, $new2)

# --- Header (docPr id="1"): BTec logo image2.jpg -> image1.jpg ---
$old3 = 'wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"'
$new3 = 'wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"'
$x = $x.Replace($old3, $new3)

$d.WordOpenXML = $x

Write-Output "done"
